$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.549.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.331.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.25%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'575.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.26%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'177.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.90%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.616"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.26%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'3.327.96"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.35%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "  -1.61%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'6.86"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.409"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'3.917.48"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.06%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'28.66"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.71%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'65.619.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "  -1.48%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.327.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.39%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'5.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.77%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'13.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.84%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'362.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.07%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'7.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.57%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'71.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.80%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.995"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.45%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.519"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.13%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "  -2.59%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.45%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.56%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'5.65"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.96%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "'22.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'6.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.25%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "  -5.48%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.45%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'159.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.56%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.847"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.72%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'26.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.45%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'1.74"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.28%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'2.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.13%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'2.671.76"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.57%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "RenderToken"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'6.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.96%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "Filecoin"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'4.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.06%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "Hedera"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'0.0663"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.70%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "OKB"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'39.62"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "Bittensor"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'333.57"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.64%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'24.26"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.72%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "  -3.55%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "  +2.58%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.965"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.65%  "
$ws.Range("E51").Style = "Normal"
